$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04872884905955236
$ws.Range("C2").Value = 0.9985669676829518

$ws.Range("B3").Value = 0.07365379909620806
$ws.Range("C3").Value = 0.9993013607077695

$ws.Range("B4").Value = 0.08899196134013831
$ws.Range("C4").Value = 0.9987951037764495
